# Update the Correspond Handoff/Handback datetimes for the
# "7f6d9555-0a25-46f2-a4db-e093a385cca0" handback entries, as part of
# regenerating the handback status report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-01-25 07:45:15"
$wsZhCn.Range("G4").Value = "2016-01-25 07:46:03"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-01-25 07:45:28"
$wsDeDe.Range("G4").Value = "2016-01-25 07:46:24"
